# Refresh the cryptos price list with the latest scrape, per the commit diff.
#
# Column D ("Price") holds plain-text numeric-looking strings such as "1.001",
# "0.01900" and "0.000007303". Assigning those verbatim to Range.Value would let
# Excel auto-detect them as numbers -- silently dropping trailing zeros and/or
# switching to scientific notation -- so each such value is written with a
# leading apostrophe (the same text-entry trick a user would type in the UI) to
# force it to stay literal text, then ClearFormats() drops the transient
# quote-prefix cell style Excel attaches for that, so the cell ends up with the
# same (default) style as before -- only its text content changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.405.74'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '1.849.87'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''233.22'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '''0.4740'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.50%  '
$ws.Range('E8').Value = '  +2.04%  '
$ws.Range('D9').Value = '''0.06316'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.82%  '
$ws.Range('D10').Value = '''17.59'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +10.26%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '''0.07469'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.59%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.815.79'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '''4.944'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = '''84.34'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.96%  '
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '30.366.16'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('D17').Value = '''243.99'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +7.65%  '
$ws.Range('D18').Value = '''1.002'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = '''12.65'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.22%  '
$ws.Range('D20').Value = '''0.000007303'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').Value = '''1.003'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = '''4.900'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.99%  '
$ws.Range('D23').Value = '''5.896'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('D24').Value = '''164.88'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '''9.075'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('E26').Value = '  +1.96%  '
$ws.Range('D27').Value = '''1.863'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.36%  '
$ws.Range('D28').Value = '''0.1031'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.46%  '
$ws.Range('D29').Value = '''1.348'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.33%  '
$ws.Range('D30').Value = '''4.031'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('D31').Value = '''3.814'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('D32').Value = '''0.04829'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('D33').Value = '''1.125'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').Value = '''0.6951'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('D36').Value = '''0.01900'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('D37').Value = '''2.677'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.60%  '
$ws.Range('D38').Value = '''1.997'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.73%  '
$ws.Range('D39').Value = '''0.8731'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.22%  '
$ws.Range('D40').Value = '''106.27'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.09%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '''0.4045'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.79%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''5.492'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('D44').Value = '''7.142'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.95%  '
$ws.Range('D45').Value = '''62.55'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.20%  '
$ws.Range('D46').Value = '''0.1196'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('E47').Value = '  +3.66%  '
$ws.Range('D48').Value = '''8.564'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.04%  '
$ws.Range('D49').Value = '''0.05514'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').Value = '''1.346'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').Value = '''0.3669'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.67%  '
